$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = "Good"
    3  = "No"
    4  = "Good"
    5  = "Good"
    6  = "Good"
    7  = "No"
    8  = "Good"
    9  = "Good"
    10 = "No"
    11 = "Good"
    12 = "No"
    13 = "No"
    14 = "Good"
    15 = "No"
    16 = "No"
    17 = "No"
    18 = "No"
    19 = "Good"
    20 = "Good"
    21 = "Good"
    22 = "Good"
    23 = "No"
    24 = "Good"
    25 = "Good"
    26 = "No"
    27 = "No"
}

foreach ($row in $values.Keys) {
    $ws.Range("I$row").Value = $values[$row]
}
